$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    if ($val -match '^-?\d+(\.\d+)?$') {
        $ws.Range($cellRef).NumberFormat = "@"
        $ws.Range($cellRef).Value = $val
        $ws.Range($cellRef).ClearFormats()
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

# --- Simple price (D column) updates ---
Set-TextValue $ws "D2" "246.15"
Set-TextValue $ws "D3" "22.68"
Set-TextValue $ws "D4" "5.278"
Set-TextValue $ws "D7" "0.8099"
Set-TextValue $ws "D8" "0.8807"
Set-TextValue $ws "D10" "0.07372"
Set-TextValue $ws "D40" "0.03912"
Set-TextValue $ws "D41" "0.006726"
Set-TextValue $ws "D44" "0.007499"
Set-TextValue $ws "D45" "0.00005626"
Set-TextValue $ws "D48" "0.1750"

# --- E column label-only updates (Best/Worst in 24h annotation changes) ---
Set-TextValue $ws "E43" "42CEJICEJI"
Set-TextValue $ws "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue $ws "E48" "47BOLOBOLO"

# --- Rows 11-26: ranking shuffled (ProBitToken moved to top), full B/C/D/E update ---
Set-TextValue $ws "B11" "ProBitToken"
Set-TextValue $ws "C11" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D11" "0.1350"
Set-TextValue $ws "E11" "10ProBitTokenPROB"
Set-TextValue $ws "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D12" "0.03043"
Set-TextValue $ws "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws "B13" "BitrueCoin"
Set-TextValue $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D13" "0.03126"
Set-TextValue $ws "E13" "12BitrueCoinBTR"
Set-TextValue $ws "B14" "BitMartToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D14" "0.09398"
Set-TextValue $ws "E14" "13BitMartTokenBMX"
Set-TextValue $ws "B15" "MCDex"
Set-TextValue $ws "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D15" "3.938"
Set-TextValue $ws "E15" "14MCDexMCB"
Set-TextValue $ws "B16" "BitForexToken"
Set-TextValue $ws "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D16" "0.001574"
Set-TextValue $ws "E16" "15BitForexTokenBF"
Set-TextValue $ws "B17" "CoinExToken"
Set-TextValue $ws "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D17" "0.04793"
Set-TextValue $ws "E17" "16CoinExTokenCET"
Set-TextValue $ws "B18" "One"
Set-TextValue $ws "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0005842"
Set-TextValue $ws "E18" "17OneONEWorstin24h"
Set-TextValue $ws "B19" "TigerCash"
Set-TextValue $ws "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D19" "0.006140"
Set-TextValue $ws "E19" "18TigerCashTCH"
Set-TextValue $ws "B20" "HotbitToken"
Set-TextValue $ws "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D20" "0.005113"
Set-TextValue $ws "E20" "19HotbitTokenHTB"
Set-TextValue $ws "B21" "BitKan"
Set-TextValue $ws "C21" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D21" "0.0009991"
Set-TextValue $ws "E21" "20BitKanKAN"
Set-TextValue $ws "B22" "NitroEx"
Set-TextValue $ws "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D22" "0.0001501"
Set-TextValue $ws "E22" "21NitroExNTX"
Set-TextValue $ws "B23" "LEO"
Set-TextValue $ws "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D23" "3.745"
Set-TextValue $ws "E23" "22LEOLEO"
Set-TextValue $ws "B24" "KuCoinToken"
Set-TextValue $ws "C24" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D24" "6.314"
Set-TextValue $ws "E24" "23KuCoinTokenKCS"
Set-TextValue $ws "B25" "BTSEToken"
Set-TextValue $ws "C25" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D25" "2.187"
Set-TextValue $ws "E25" "24BTSETokenBTSE"
Set-TextValue $ws "B26" "BitpandaEcosystemToken"
Set-TextValue $ws "C26" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D26" "0.3280"
Set-TextValue $ws "E26" "25BitpandaEcosystemTokenBEST"

Write-Host "Edit complete"
